# Update column F ("report count" style metric) values across all four
# sheets of the 北京-漫展信息 workbook, matching the upstream regeneration
# at commit 456a3b4 ("Update gh-pages to output generated at 456a3b4").
#
# Column F holds a numeric counter for each listing row; every other
# column (A-E, G-I) is untouched by this refresh.

$wb = $excel.ActiveWorkbook

# Column F is index 6 (A=1 ... F=6)
$col = 6

# Sheet "展览" (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$updates = @{
    3  = 135
    6  = 261
    7  = 13090
    8  = 61
    10 = 266
    11 = 3270
    13 = 6582
    15 = 49
    16 = 3470
    17 = 167
    19 = 28
    21 = 67
    23 = 51
    24 = 3641
    25 = 100
    27 = 3075
    28 = 419
    29 = 1902
    30 = 104
    31 = 225
    32 = 6759
    34 = 1254
    35 = 2009
    36 = 1298
    37 = 106
    38 = 1057
    40 = 217
    43 = 1145
    44 = 142
    45 = 1215
    46 = 1800
    48 = 162
    49 = 1176
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, $col).Value = $updates[$row]
}

# Sheet "演出" (Performances)
$ws = $wb.Worksheets.Item("演出")
$updates = @{
    8  = 127
    14 = 103
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, $col).Value = $updates[$row]
}

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$updates = @{
    2 = 451
    3 = 620
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, $col).Value = $updates[$row]
}

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$updates = @{
    3  = 135
    6  = 451
    7  = 620
    8  = 261
    9  = 13090
    10 = 61
    13 = 266
    14 = 3270
    16 = 3470
    17 = 167
    18 = 28
    20 = 67
    23 = 51
    24 = 3641
    27 = 3075
    28 = 3076
    29 = 419
    30 = 1902
    31 = 104
    32 = 225
    33 = 6759
    34 = 103
    36 = 1255
    37 = 2009
    39 = 1298
    40 = 106
    41 = 1057
    42 = 217
    45 = 1215
    47 = 1800
    50 = 162
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, $col).Value = $updates[$row]
}
